# Auto-generated script applying the Carbuncle_Profits leve-price/profit data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2201.9092
$ws.Range("I6").Value = 333.16666
$ws.Range("J6").Value = 4444.4
$ws.Range("K6").Value = 999.4999799999999
$ws.Range("L6").Value = 13333.2
$ws.Range("M6").Value = -887.4999799999999
$ws.Range("N6").Value = -13557.2

$ws.Range("H53").Value = 231.88889
$ws.Range("I53").Value = 230.6923
$ws.Range("J53").Value = 235
$ws.Range("K53").Value = 230.6923
$ws.Range("L53").Value = 235
$ws.Range("M53").Value = 406.3077
$ws.Range("N53").Value = -1509

$ws.Range("H86").Value = 46533.668
$ws.Range("I86").Value = 2083.1667
$ws.Range("J86").Value = 135434.67
$ws.Range("K86").Value = 2083.1667
$ws.Range("L86").Value = 135434.67
$ws.Range("M86").Value = -960.1667000000002
$ws.Range("N86").Value = -137680.67

$ws.Range("H89").Value = 46533.668
$ws.Range("I89").Value = 2083.1667
$ws.Range("J89").Value = 135434.67
$ws.Range("K89").Value = 10415.8335
$ws.Range("L89").Value = 677173.3500000001
$ws.Range("M89").Value = -4799.833500000001
$ws.Range("N89").Value = -688405.3500000001

$ws.Range("H112").Value = 1317.0588
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1372.6666
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 4117.9998
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -6333.9998

$ws.Range("H129").Value = 871.9268
$ws.Range("J129").Value = 876.7625
$ws.Range("L129").Value = 2630.2875
$ws.Range("N129").Value = -12630.2875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 500000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H122").Value = 4460
$ws.Range("I122").Value = 3950
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 11850
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -9400
$ws.Range("N122").Value = -19300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2530.1428
$ws.Range("I86").Value = 3142.2
$ws.Range("K86").Value = 3142.2
$ws.Range("M86").Value = -2019.2

$ws.Range("H89").Value = 2530.1428
$ws.Range("I89").Value = 3142.2
$ws.Range("K89").Value = 15711
$ws.Range("M89").Value = -10095

$ws.Range("H105").Value = 3454.9
$ws.Range("I105").Value = 2727.6667
$ws.Range("K105").Value = 2727.6667
$ws.Range("M105").Value = -980.6667000000002

$ws.Range("H134").Value = 2227.2058
$ws.Range("I134").Value = 1065.0869
$ws.Range("J134").Value = 4657.091
$ws.Range("K134").Value = 3195.2607
$ws.Range("L134").Value = 13971.273
$ws.Range("M134").Value = -660.2606999999998
$ws.Range("N134").Value = -19041.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3752
$ws.Range("I132").Value = 4331.6
$ws.Range("J132").Value = 3529.077
$ws.Range("K132").Value = 12994.8
$ws.Range("L132").Value = 10587.231
$ws.Range("M132").Value = -10464.8
$ws.Range("N132").Value = -15647.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1026.6666
$ws.Range("I59").Value = 700
$ws.Range("J59").Value = 1680
$ws.Range("K59").Value = 2100
$ws.Range("L59").Value = 5040
$ws.Range("M59").Value = -1560
$ws.Range("N59").Value = -6120

$ws.Range("H98").Value = 1099.8334
$ws.Range("J98").Value = 1099.8334
$ws.Range("L98").Value = 3299.5002
$ws.Range("N98").Value = -6295.5002

$ws.Range("H117").Value = 1764.5
$ws.Range("I117").Value = 529
$ws.Range("K117").Value = 1587
$ws.Range("M117").Value = 1855

$ws.Range("H129").Value = 1881.8
$ws.Range("I129").Value = 1166
$ws.Range("J129").Value = 2120.4
$ws.Range("K129").Value = 3498
$ws.Range("L129").Value = 6361.200000000001
$ws.Range("M129").Value = 1502
$ws.Range("N129").Value = -16361.2

$ws.Range("H131").Value = 988.71
$ws.Range("I131").Value = 722
$ws.Range("J131").Value = 1002.7474
$ws.Range("K131").Value = 2166
$ws.Range("L131").Value = 3008.2422
$ws.Range("M131").Value = 2874
$ws.Range("N131").Value = -13088.2422

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6500
$ws.Range("I80").Value = 6500
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 6500
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -5502
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 6500
$ws.Range("I83").Value = 6500
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 32500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -27508
$ws.Range("N83").ClearContents()

$ws.Range("H102").Value = 1460
$ws.Range("I102").Value = 1433.3334
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1433.3334
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 188.6666
$ws.Range("N102").Value = -4744

$ws.Range("H122").Value = 2328.818
$ws.Range("I122").Value = 1939.9565
$ws.Range("J122").Value = 3223.2
$ws.Range("K122").Value = 5819.8695
$ws.Range("L122").Value = 9669.599999999999
$ws.Range("M122").Value = -3369.8695
$ws.Range("N122").Value = -14569.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 78431.695
$ws.Range("I7").Value = 112123.555
$ws.Range("J7").Value = 2625
$ws.Range("K7").Value = 112123.555
$ws.Range("L7").Value = 2625
$ws.Range("M7").Value = -112011.555
$ws.Range("N7").Value = -2849

$ws.Range("H16").Value = 1464.2354
$ws.Range("I16").Value = 823.75
$ws.Range("J16").Value = 2033.5555
$ws.Range("K16").Value = 823.75
$ws.Range("L16").Value = 2033.5555
$ws.Range("M16").Value = -653.75
$ws.Range("N16").Value = -2373.5555

$ws.Range("H126").Value = 78431.695
$ws.Range("I126").Value = 112123.555
$ws.Range("J126").Value = 2625
$ws.Range("K126").Value = 336370.665
$ws.Range("L126").Value = 7875
$ws.Range("M126").Value = -333900.665
$ws.Range("N126").Value = -12815

$ws.Range("H136").Value = 27784624
$ws.Range("I136").Value = 5088.3335
$ws.Range("K136").Value = 15265.0005
$ws.Range("M136").Value = -12715.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 73614.5
$ws.Range("I122").Value = 73614.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 220843.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -218393.5
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 146400.86
$ws.Range("I126").Value = 252500.25
$ws.Range("J126").Value = 4935
$ws.Range("K126").Value = 757500.75
$ws.Range("L126").Value = 14805
$ws.Range("M126").Value = -755030.75
$ws.Range("N126").Value = -19745

